$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2115384615384615
$ws.Range("C2").Value = 0.5192307692307693
$ws.Range("J2").Value = 0.01538461538461539
$ws.Range("P2").Value = 0.1538461538461539
$ws.Range("S2").Value = 0.1
$ws.Range("B3").Value = 0.01379310344827586
$ws.Range("C3").Value = 0.04137931034482759
$ws.Range("J3").Value = 0.006896551724137931
$ws.Range("P3").Value = 0.7103448275862069
$ws.Range("S3").Value = 0.2275862068965517
$ws.Range("J4").Value = 0.06818181818181818
$ws.Range("P4").Value = 0.5681818181818182
$ws.Range("S4").Value = 0.3636363636363636
$ws.Range("B6").Value = 0.05627705627705628
$ws.Range("D6").Value = 0.004329004329004329
$ws.Range("F6").Value = 0.04329004329004329
$ws.Range("J6").Value = 0.3766233766233766
$ws.Range("O6").Value = 0.02597402597402598
$ws.Range("Q6").Value = 0.1341991341991342
$ws.Range("R6").Value = 0.04329004329004329
$ws.Range("S6").Value = 0.316017316017316
$ws.Range("B7").Value = 0.06912442396313365
$ws.Range("D7").Value = 0.0184331797235023
$ws.Range("F7").Value = 0.04147465437788019
$ws.Range("J7").Value = 0.1290322580645161
$ws.Range("O7").Value = 0.02304147465437788
$ws.Range("Q7").Value = 0.1981566820276498
$ws.Range("R7").Value = 0.06451612903225806
$ws.Range("S7").Value = 0.4562211981566821
$ws.Range("B8").Value = 0.07468879668049792
$ws.Range("D8").Value = 0.02697095435684647
$ws.Range("E8").Value = 0.002074688796680498
$ws.Range("F8").Value = 0.04771784232365145
$ws.Range("J8").Value = 0.06431535269709543
$ws.Range("O8").Value = 0.02282157676348548
$ws.Range("Q8").Value = 0.2116182572614108
$ws.Range("R8").Value = 0.0954356846473029
$ws.Range("S8").Value = 0.454356846473029
$ws.Range("B9").Value = 0.07471264367816093
$ws.Range("D9").Value = 0.04022988505747126
$ws.Range("F9").Value = 0.06321839080459771
$ws.Range("J9").Value = 0.07471264367816093
$ws.Range("O9").Value = 0.005747126436781609
$ws.Range("Q9").Value = 0.1436781609195402
$ws.Range("R9").Value = 0.08620689655172414
$ws.Range("S9").Value = 0.5114942528735632
$ws.Range("B10").Value = 0.1035714285714286
$ws.Range("D10").Value = 0.01696428571428571
$ws.Range("E10").Value = 0.001785714285714286
$ws.Range("F10").Value = 0.08125
$ws.Range("J10").Value = 0.09285714285714286
$ws.Range("O10").Value = 0.01964285714285714
$ws.Range("Q10").Value = 0.2107142857142857
$ws.Range("R10").Value = 0.07410714285714286
$ws.Range("S10").Value = 0.3991071428571429
$ws.Range("G11").Value = 0.1268011527377522
$ws.Range("J11").Value = 0.08645533141210375
$ws.Range("K11").Value = 0.1959654178674352
$ws.Range("L11").Value = 0.5619596541786743
$ws.Range("S11").Value = 0.02881844380403458
$ws.Range("G12").Value = 0.7424242424242424
$ws.Range("J12").Value = 0.202020202020202
$ws.Range("K12").Value = 0.0101010101010101
$ws.Range("L12").Value = 0.0101010101010101
$ws.Range("S12").Value = 0.03535353535353535
$ws.Range("G13").Value = 0.7391304347826086
$ws.Range("J13").Value = 0.2173913043478261
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("F15").Value = 0.0186046511627907
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.05581395348837209
$ws.Range("J15").Value = 0.3395348837209302
$ws.Range("K15").Value = 0.06046511627906977
$ws.Range("O15").Value = 0.04651162790697674
$ws.Range("S15").Value = 0.2790697674418605
$ws.Range("F16").Value = 0.04216867469879518
$ws.Range("H16").Value = 0.2168674698795181
$ws.Range("I16").Value = 0.0963855421686747
$ws.Range("J16").Value = 0.2951807228915663
$ws.Range("K16").Value = 0.108433734939759
$ws.Range("M16").Value = 0.03614457831325301
$ws.Range("O16").Value = 0.03012048192771084
$ws.Range("S16").Value = 0.1746987951807229
$ws.Range("F17").Value = 0.03023255813953488
$ws.Range("H17").Value = 0.2046511627906977
$ws.Range("I17").Value = 0.07209302325581396
$ws.Range("J17").Value = 0.4069767441860465
$ws.Range("K17").Value = 0.09302325581395349
$ws.Range("M17").Value = 0.01627906976744186
$ws.Range("O17").Value = 0.05813953488372093
$ws.Range("S17").Value = 0.1186046511627907
$ws.Range("F18").Value = 0.02994011976047904
$ws.Range("H18").Value = 0.155688622754491
$ws.Range("I18").Value = 0.1137724550898204
$ws.Range("J18").Value = 0.3772455089820359
$ws.Range("K18").Value = 0.1017964071856287
$ws.Range("M18").Value = 0.01197604790419162
$ws.Range("O18").Value = 0.1017964071856287
$ws.Range("S18").Value = 0.1077844311377246
$ws.Range("F19").Value = 0.02828746177370031
$ws.Range("H19").Value = 0.2285932721712538
$ws.Range("I19").Value = 0.07415902140672782
$ws.Range("J19").Value = 0.331039755351682
$ws.Range("K19").Value = 0.1391437308868501
$ws.Range("M19").Value = 0.02446483180428135
$ws.Range("O19").Value = 0.06345565749235474
$ws.Range("S19").Value = 0.1108562691131498
